$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are treated as Text so numeric-looking values
# (e.g. '303.57', '43.105.21') are stored as strings, matching the original inlineStr format.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.105.21'
$ws.Range("E2").Value = '  -0.16%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.369.43'
$ws.Range("E3").Value = '  +1.18%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.57'
$ws.Range("E5").Value = '  +0.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.32'
$ws.Range("E6").Value = '  +1.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.504'
$ws.Range("E7").Value = '  -0.25%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.482'
$ws.Range("E9").Value = '  -2.82%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.40'
$ws.Range("E10").Value = '  +0.75%  '

$ws.Range("E11").Value = '  +3.65%  '

$ws.Range("E12").Value = '  +0.47%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.34'
$ws.Range("E13").Value = '  -1.96%  '

$ws.Range("E14").Value = '  +0.42%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.736.24'
$ws.Range("E15").Value = '  +1.18%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.401.59'
$ws.Range("E16").Value = '  +0.45%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.802'
$ws.Range("E17").Value = '  +0.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.132.51'
$ws.Range("E18").Value = '  +0.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.99'
$ws.Range("E19").Value = '  -1.77%  '

$ws.Range("E20").Value = '  +1.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0889'
$ws.Range("E21").Value = '  -0.25%  '

$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.69'
$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("E24").Value = '  +0.89%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.44'
$ws.Range("E25").Value = '  +1.12%  '

$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.47'
$ws.Range("E27").Value = '  -0.71%  '

$ws.Range("E28").Value = '  +0.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.36'
$ws.Range("E29").Value = '  +1.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.90'
$ws.Range("E30").Value = '  +1.30%  '

$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("E32").Value = '  +0.67%  '

$ws.Range("E33").Value = '  +11.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.83'
$ws.Range("E34").Value = '  +3.14%  '

$ws.Range("E35").Value = '  +0.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '128.28'
$ws.Range("E36").Value = '  +13.14%  '

$ws.Range("E37").Value = '  -0.46%  '

$ws.Range("E40").Value = '  -3.03%  '

$ws.Range("E41").Value = '  -0.73%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.11'
$ws.Range("E42").Value = '  -4.91%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.925.85'
$ws.Range("E43").Value = '  -0.65%  '

$ws.Range("E44").Value = '  -1.42%  '

$ws.Range("E45").Value = '  +1.85%  '

$ws.Range("E46").Value = '  +1.32%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.27'
$ws.Range("E47").Value = '  -7.48%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.596.65'
$ws.Range("E48").Value = '  +1.04%  '

$ws.Range("E49").Value = '  +3.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '71.72'
$ws.Range("E50").Value = '  -0.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.75'
$ws.Range("E51").Value = '  -2.73%  '

# Swap rows 38 and 39 (RenderToken / LidoDAOToken reorder) with updated values
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.35'
$ws.Range("E38").Value = '  -0.56%  '

$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.86'
$ws.Range("E39").Value = '  +3.84%  '

